# Project outline edit: insert "grenades, " into the "Power-ups include..." bullet,
# right before "invincibility, etc.", matching the structure produced by Word when a
# user types a short insertion in the middle of an existing sentence (the new text
# becomes its own run, and the document's _GoBack bookmark relocates to the point of
# the most recent edit).

$d = $word.ActiveDocument
$insertText = "grenades, "

# 1. Locate the insertion point: right before "invincibility, etc." in the
#    "Power-ups include ..." bullet.
$targetRng = $d.Content.Duplicate
$foundTarget = $targetRng.Find.Execute("invincibility, etc.", $true, $false, $false, $false,
                        $false, $true, 1, $false, "", 0)
if (-not $foundTarget) {
    Write-Host "ERROR: could not find target phrase 'invincibility, etc.'"
}
$targetStart = $targetRng.Start

# 2. Stage the new text in a harmless scratch spot (right after the document title)
#    so it picks up plain/default run formatting, then we can Copy it as an isolated
#    run and Paste it at the real target -- this avoids merging it into the
#    surrounding run's <w:r> element.
$titleRng = $d.Content.Duplicate
$foundTitle = $titleRng.Find.Execute("Monster Maze Project", $true, $false, $false, $false,
                        $false, $true, 1, $false, "", 0)
if (-not $foundTitle) {
    Write-Host "ERROR: could not find document title to stage scratch text"
}
$scratchStart = $titleRng.End
$scratchInsertionPoint = $d.Range($scratchStart, $scratchStart)
$scratchInsertionPoint.InsertAfter($insertText)
$scratchEnd = $scratchStart + $insertText.Length

# The scratch text was inserted earlier in the document than the target, so the
# target position shifts forward by the length of the inserted text.
$newTargetStart = $targetStart + $insertText.Length

# 3. Copy the scratch run and paste it at the real target location.
$scratchRng = $d.Range($scratchStart, $scratchEnd)
$scratchRng.Copy()

$pasteRng = $d.Range($newTargetStart, $newTargetStart)
$pasteRng.Paste()

# 4. Remove the scratch text now that it has been duplicated into place.
$scratchRngAgain = $d.Range($scratchStart, $scratchEnd)
$scratchRngAgain.Delete()

# 5. Move the _GoBack bookmark so it sits between the newly-inserted "grenades, "
#    run and "invincibility, etc." (adding a bookmark with an existing name moves it,
#    since bookmark names are unique within the document).
$bmRng = $d.Range($newTargetStart, $newTargetStart)
$d.Bookmarks.Add("_GoBack", $bmRng)

Write-Host "Edit complete."
